$d = $word.ActiveDocument

# 1. Fix schedule time
$d.Content.Find.Execute("19:30-20:25 Cleaning and Transforming", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "19:30-20:15 Cleaning and Transforming", 2)

# 2. Shorten regex description
$d.Content.Find.Execute(": exclude match regular expressions", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ": exclude match", 2)

# 3. Replace windows path with relative path
$d.Content.Find.Execute("c:\users\DSTraining\desktop\119\text\", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ".../libcarp-wk2-data/text", 2)

# 4. Extend sed sentence
$d.Content.Find.Execute("command allows you to edit files directly.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "command allows you to edit files directly. This can be used to remove all the header and footer information that Project Gutenberg add before and after a text.", 2)

# 5. Replace Notepad++ mentions with "a text editor"
$d.Content.Find.Execute("in Notepad++. Note how the text has been transformed ready for analysis.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "in a text editor. Note how the text has been transformed ready for analysis.", 2)

$d.Content.Find.Execute("Open the file in Notepad++ and after scrolling past some blank space", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Open the file in a text editor and after scrolling past some blank space", 2)
